$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) stays text, matching original inlineStr formatting
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "23.419.35"
$ws.Range("E2").Value = "  +0.72%  "
$ws.Range("D3").Value = "1.640.29"
$ws.Range("E3").Value = "  +2.32%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "1.001"
$ws.Range("E5").Value = "  +0.02%  "
$ws.Range("D6").Value = "305.02"
$ws.Range("E6").Value = "  +0.45%  "
$ws.Range("D7").Value = "0.3734"
$ws.Range("E7").Value = "  -0.97%  "
$ws.Range("D8").Value = "51.98"
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "0.3626"
$ws.Range("E9").Value = "  -0.57%  "
$ws.Range("D10").Value = "1.243"
$ws.Range("E10").Value = "  -2.85%  "
$ws.Range("D11").Value = "0.08098"
$ws.Range("E11").Value = "  -0.55%  "
$ws.Range("E12").Value = "  +0.03%  "
$ws.Range("D13").Value = "22.71"
$ws.Range("E13").Value = "  -0.41%  "
$ws.Range("D14").Value = "6.586"
$ws.Range("E14").Value = "  -0.16%  "
$ws.Range("D15").Value = "0.00001267"
$ws.Range("E15").Value = "  +1.51%  "
$ws.Range("D16").Value = "7.258"
$ws.Range("E16").Value = "  -2.33%  "
$ws.Range("D17").Value = "1.634.25"
$ws.Range("E17").Value = "  +1.94%  "
$ws.Range("D18").Value = "94.60"
$ws.Range("E18").Value = "  +0.57%  "
$ws.Range("D19").Value = "0.06878"
$ws.Range("E19").Value = "  -0.73%  "
$ws.Range("D20").Value = "18.09"
$ws.Range("E20").Value = "  -0.49%  "
$ws.Range("D21").Value = "6.505"
$ws.Range("E21").Value = "  -0.55%  "
$ws.Range("D22").Value = "1.001"
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("D23").Value = "23.427.78"
$ws.Range("E23").Value = "  +0.73%  "
$ws.Range("E24").Value = "  -1.36%  "
$ws.Range("D25").Value = "3.097"
$ws.Range("E25").Value = "  +0.43%  "
$ws.Range("D26").Value = "2.399"
$ws.Range("E26").Value = "  +0.74%  "
$ws.Range("D27").Value = "21.17"
$ws.Range("E27").Value = "  -0.20%  "
$ws.Range("D28").Value = "150.66"
$ws.Range("E28").Value = "  +0.71%  "
$ws.Range("D29").Value = "5.278"
$ws.Range("E29").Value = "  +0.41%  "
$ws.Range("D30").Value = "135.94"
$ws.Range("E30").Value = "  +0.95%  "
$ws.Range("D31").Value = "2.288"
$ws.Range("E31").Value = "  -4.44%  "
$ws.Range("D32").Value = "1.814.93"
$ws.Range("E32").Value = "  +1.84%  "
$ws.Range("D33").Value = "6.837"
$ws.Range("E33").Value = "  +1.42%  "
$ws.Range("D34").Value = "0.9495"
$ws.Range("E34").Value = "  -1.58%  "
$ws.Range("D35").Value = "0.02796"
$ws.Range("E35").Value = "  +1.41%  "
$ws.Range("D36").Value = "10.42"
$ws.Range("E36").Value = "  +0.69%  "
$ws.Range("D37").Value = "0.2523"
$ws.Range("E37").Value = "  -0.52%  "
$ws.Range("B38").Value = "InternetComputer(DFINITY)"
$ws.Range("C38").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D38").Value = "6.158"
$ws.Range("E38").Value = "  +0.91%  "
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").Value = "0.07238"
$ws.Range("E39").Value = "  -3.63%  "
$ws.Range("D40").Value = "0.08753"
$ws.Range("E40").Value = "  -0.79%  "
$ws.Range("D41").Value = "1.368"
$ws.Range("E41").Value = "  -1.14%  "
$ws.Range("D42").Value = "0.7033"
$ws.Range("E42").Value = "  -1.19%  "
$ws.Range("D43").Value = "12.45"
$ws.Range("E43").Value = "  -0.17%  "
$ws.Range("D44").Value = "16.03"
$ws.Range("E44").Value = "  +3.09%  "
$ws.Range("D45").Value = "0.6498"
$ws.Range("E45").Value = "  -1.01%  "
$ws.Range("D46").Value = "2.322"
$ws.Range("E46").Value = "  +0.14%  "
$ws.Range("D47").Value = "0.9995"
$ws.Range("D48").Value = "4.005"
$ws.Range("E48").Value = "  -0.21%  "
$ws.Range("D49").Value = "0.07969"
$ws.Range("E49").Value = "  +0.31%  "
$ws.Range("D50").Value = "128.28"
$ws.Range("E50").Value = "  -3.23%  "
$ws.Range("D51").Value = "1.201"
$ws.Range("E51").Value = "  -0.48%  "
